$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C9").Value = 8696
$ws.Range("C10:C11").Value = 8573
$ws.Range("C12").Value = 8241
$ws.Range("C13:C52").Value = 8183
$ws.Range("C53:C76").Value = 7691
$ws.Range("C77:C252").Value = 7622
